# Fruta / hortaliza, semanal
# The underlying dataset rows (2-21) got re-shuffled: each target row now
# holds the full record (columns A:R) that used to live in a different
# source row. This reproduces that reshuffle by snapshotting every row's
# values first (so reads never see already-overwritten data) and then
# writing each target row from its mapped source snapshot.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 2
$lastRow = 21
$lastCol = 18   # column R

# Map: target row -> source row (where the data used to live before the edit)
$rowMap = @{
    2  = 20
    3  = 3
    4  = 19
    5  = 11
    6  = 5
    7  = 15
    8  = 6
    9  = 4
    10 = 10
    11 = 13
    12 = 16
    13 = 7
    14 = 21
    15 = 2
    16 = 17
    17 = 14
    18 = 9
    19 = 12
    20 = 18
    21 = 8
}

# Snapshot every row's current values (A:R) before writing anything back.
$snapshot = @{}
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $rowVals = @()
    for ($c = 1; $c -le $lastCol; $c++) {
        $rowVals += ,$ws.Cells.Item($r, $c).Value2
    }
    $snapshot[$r] = $rowVals
}

# Write each target row using the values captured from its source row.
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $srcRow = $rowMap[$r]
    $srcVals = $snapshot[$srcRow]
    for ($c = 1; $c -le $lastCol; $c++) {
        $ws.Cells.Item($r, $c).Value2 = $srcVals[$c - 1]
    }
}
